$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0276967930029155
$ws.Range("C2").Value = 0.188775510204082
$ws.Range("D2").Value = 0.010932944606414
$ws.Range("E2").Value = 0.737609329446064
$ws.Range("F2").Value = 0.239067055393586
$ws.Range("G2").Value = 0.644314868804665
$ws.Range("H2").Value = 0.302478134110787
$ws.Range("I2").Value = 0.716472303206997
$ws.Range("J2").Value = 0.00947521865889213
$ws.Range("K2").Value = 0.749271137026239
$ws.Range("L2").Value = 0.754373177842566
$ws.Range("M2").Value = 0.0196793002915452
$ws.Range("N2").Value = 0.209183673469388
$ws.Range("O2").Value = 0.358600583090379
$ws.Range("P2").Value = 0.00145772594752187
$ws.Range("Q2").Value = 0.0043731778425656
$ws.Range("R2").Value = 0.96064139941691
$ws.Range("S2").Value = 0.123906705539359
$ws.Range("T2").Value = 0.262390670553936
$ws.Range("U2").Value = 0.0532069970845481
$ws.Range("V2").Value = 0.0291545189504373
$ws.Range("W2").Value = 0.0247813411078717
$ws.Range("X2").Value = 0.010932944606414
$ws.Range("B3").Value = 0.933673469387755
$ws.Range("C3").Value = 0.0284256559766764
$ws.Range("D3").Value = 0.0160349854227405
$ws.Range("E3").Value = 0.0153061224489796
$ws.Range("F3").Value = 0.0313411078717201
$ws.Range("G3").Value = 0.00801749271137026
$ws.Range("H3").Value = 0.00947521865889213
$ws.Range("I3").Value = 0.0116618075801749
$ws.Range("J3").Value = 0.236880466472303
$ws.Range("K3").Value = 0.0167638483965015
$ws.Range("L3").Value = 0.0116618075801749
$ws.Range("M3").Value = 0.176384839650146
$ws.Range("N3").Value = 0.0102040816326531
$ws.Range("O3").Value = 0.00364431486880466
$ws.Range("P3").Value = 0.0641399416909621
$ws.Range("Q3").Value = 0.0116618075801749
$ws.Range("R3").Value = 0.00291545189504373
$ws.Range("S3").Value = 0.0349854227405248
$ws.Range("T3").Value = 0.0145772594752187
$ws.Range("U3").Value = 0.0204081632653061
$ws.Range("V3").Value = 0.0087463556851312
$ws.Range("W3").Value = 0.0087463556851312
$ws.Range("X3").Value = 0.0196793002915452
$ws.Range("B4").Value = 0.0189504373177843
$ws.Range("C4").Value = 0.38265306122449
$ws.Range("D4").Value = 0.744897959183674
$ws.Range("E4").Value = 0.223760932944606
$ws.Range("F4").Value = 0.724489795918367
$ws.Range("G4").Value = 0.32798833819242
$ws.Range("H4").Value = 0.470845481049563
$ws.Range("I4").Value = 0.260204081632653
$ws.Range("J4").Value = 0.750728862973761
$ws.Range("K4").Value = 0.00801749271137026
$ws.Range("L4").Value = 0.0240524781341108
$ws.Range("M4").Value = 0.766034985422741
$ws.Range("N4").Value = 0.0473760932944606
$ws.Range("O4").Value = 0.637026239067055
$ws.Range("P4").Value = 0.217930029154519
$ws.Range("Q4").Value = 0.00291545189504373
$ws.Range("R4").Value = 0.0276967930029155
$ws.Range("S4").Value = 0.837463556851312
$ws.Range("T4").Value = 0.0189504373177843
$ws.Range("U4").Value = 0.91399416909621
$ws.Range("V4").Value = 0.0145772594752187
$ws.Range("W4").Value = 0.957725947521866
$ws.Range("X4").Value = 0.739067055393586
$ws.Range("B5").Value = 0.0196793002915452
$ws.Range("C5").Value = 0.399416909620991
$ws.Range("D5").Value = 0.228134110787172
$ws.Range("E5").Value = 0.0233236151603499
$ws.Range("F5").Value = 0.0043731778425656
$ws.Range("G5").Value = 0.0196793002915452
$ws.Range("H5").Value = 0.216472303206997
$ws.Range("I5").Value = 0.0116618075801749
$ws.Range("J5").Value = 0.00291545189504373
$ws.Range("K5").Value = 0.225947521865889
$ws.Range("L5").Value = 0.209912536443149
$ws.Range("M5").Value = 0.0379008746355685
$ws.Range("N5").Value = 0.733236151603499
$ws.Range("P5").Value = 0.715014577259475
$ws.Range("Q5").Value = 0.981049562682216
$ws.Range("R5").Value = 0.00801749271137026
$ws.Range("S5").Value = 0.00291545189504373
$ws.Range("T5").Value = 0.704081632653061
$ws.Range("U5").Value = 0.0123906705539359
$ws.Range("V5").Value = 0.947521865889213
$ws.Range("W5").Value = 0.00801749271137026
$ws.Range("X5").Value = 0.229591836734694
